$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C header (13-01-2023), matching B1's format ---
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: 1810 Renta variable ---
$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 1263.59
$ws.Range("C2").Value = 1519.64

# --- Row 3: Alpha Acciones ---
$ws.Range("A3").Value = "Alpha Acciones"
$ws.Range("B3").Value = 54394.51
$ws.Range("C3").Value = 54407.67

# --- Row 4: Alpha Mega ---
$ws.Range("A4").Value = "Alpha Mega"
$ws.Range("B4").Value = 103608.94
$ws.Range("C4").Value = 103625.92

# --- Row 5: Fima Acciones ---
$ws.Range("A5").Value = "Fima Acciones"
$ws.Range("B5").Value = 9630.27
$ws.Range("C5").Value = 9805.98

# --- Row 6: Fima PB Acciones ---
$ws.Range("A6").Value = "Fima PB Acciones"
$ws.Range("B6").Value = 4295.94
$ws.Range("C6").Value = 4026.55

# --- Row 7: HF Acciones Argentinas ---
$ws.Range("A7").Value = "HF Acciones Argentinas"
$ws.Range("B7").Value = 408.3
$ws.Range("C7").Value = 445.02

# --- Row 8: avg ---
$ws.Range("A8").Value = "avg"
$ws.Range("B8").Value = 28933.59
$ws.Range("C8").Value = 28971.8

# --- Row 9: total ---
$ws.Range("A9").Value = "total"
$ws.Range("B9").Value = 173601.55
$ws.Range("C9").Value = 173830.78
